$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: new attendance date column F (17-05-2025), matching the
# bold/centered/bordered style already used by the other header cells ---
$ws.Range("F1").Value = "17-05-2025"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)   # xlPasteFormats - copy format only, no value
$ws.Range("F1").Value = "17-05-2025"

# --- Updated attendance percentages for row 2 and row 4 ---
$ws.Range("C2").Value = 33.33333333333333
$ws.Range("C4").Value = 33.33333333333333

# --- New attendance marks in column F for the existing students ---
$ws.Range("F2").Value = "-"
$ws.Range("F3").Value = "✓"
$ws.Range("F4").Value = "-"

# --- New student row 5 ---
# A5 should stay text ("156"), like the other roll numbers in column A.
$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "156"
$ws.Range("A2").Copy()
$ws.Range("A5").PasteSpecial(-4122)   # strip the temporary text-format style again

$ws.Range("B5").Value = "P. Diddy"
$ws.Range("C5").Value = 0

# D5/E5/F5 have no attendance recorded yet - leave them blank, but still
# present on the sheet (a tiny formatting touch keeps them from being
# dropped entirely when the row is saved).
$ws.Range("D5").Font.Bold = $false
$ws.Range("E5").Font.Bold = $false
$ws.Range("F5").Font.Bold = $false
